$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data (row 9): task name and hours
$ws.Range("B9").Value = "Login with JWT and logout"
$ws.Range("C9").Value = "17,00-"

# Move the active selection to the newly added cell, matching the saved view state
$null = $ws.Range("C9").Select()
